$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 600
$ws.Range("I6").Value = 675
$ws.Range("K6").Value = 2025
$ws.Range("M6").Value = -1913
$ws.Range("H11").Value = 124.4
$ws.Range("I11").Value = 124.4
$ws.Range("K11").Value = 124.4
$ws.Range("M11").Value = 15.59999999999999
$ws.Range("H40").Value = 3710.7144
$ws.Range("I40").Value = 3358.3333
$ws.Range("J40").Value = 3975
$ws.Range("K40").Value = 3358.3333
$ws.Range("L40").Value = 3975
$ws.Range("M40").Value = -3183.3333
$ws.Range("N40").Value = -4325
$ws.Range("H95").Value = 36674.668
$ws.Range("J95").Value = 36674.668
$ws.Range("L95").Value = 36674.668
$ws.Range("N95").Value = -42166.668
$ws.Range("H137").Value = 2153.5405
$ws.Range("I137").Value = 1946.9286
$ws.Range("J137").Value = 2796.3333
$ws.Range("K137").Value = 5840.7858
$ws.Range("L137").Value = 8388.999899999999
$ws.Range("M137").Value = -3290.7858
$ws.Range("N137").Value = -13488.9999
$ws.Range("H138").Value = 2744.147
$ws.Range("I138").Value = 1550.5
$ws.Range("J138").Value = 3395.2273
$ws.Range("K138").Value = 4651.5
$ws.Range("L138").Value = 10185.6819
$ws.Range("M138").Value = 488.5
$ws.Range("N138").Value = -20465.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3126.3
$ws.Range("I45").Value = 2170
$ws.Range("K45").Value = 2170
$ws.Range("M45").Value = -1793
$ws.Range("H61").Value = 26112.5
$ws.Range("I61").Value = 26112.5
$ws.Range("K61").Value = 26112.5
$ws.Range("M61").Value = -25900.5
$ws.Range("H97").Value = 3535.25
$ws.Range("I97").Value = 1425.1428
$ws.Range("J97").Value = 6489.4
$ws.Range("K97").Value = 1425.1428
$ws.Range("L97").Value = 6489.4
$ws.Range("M97").Value = -929.1428000000001
$ws.Range("N97").Value = -7481.4
$ws.Range("H102").Value = 3349.2856
$ws.Range("I102").Value = 3597.3076
$ws.Range("J102").Value = 125
$ws.Range("K102").Value = 3597.3076
$ws.Range("L102").Value = 125
$ws.Range("M102").Value = -1975.3076
$ws.Range("N102").Value = -3369
$ws.Range("H122").Value = 3840.7666
$ws.Range("I122").Value = 3408.92
$ws.Range("K122").Value = 10226.76
$ws.Range("M122").Value = -7776.76
$ws.Range("H136").Value = 26112.5
$ws.Range("I136").Value = 26112.5
$ws.Range("K136").Value = 78337.5
$ws.Range("M136").Value = -75787.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 14168.25
$ws.Range("I107").Value = 14438.25
$ws.Range("J107").Value = 13628.25
$ws.Range("K107").Value = 14438.25
$ws.Range("L107").Value = 13628.25
$ws.Range("M107").Value = -12518.25
$ws.Range("N107").Value = -17468.25
$ws.Range("H134").Value = 2490.1428
$ws.Range("I134").Value = 2063.3635
$ws.Range("J134").Value = 2959.6
$ws.Range("K134").Value = 6190.0905
$ws.Range("L134").Value = 8878.799999999999
$ws.Range("M134").Value = -3655.0905
$ws.Range("N134").Value = -13948.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2188.2273
$ws.Range("J31").Value = 4999.6665
$ws.Range("L31").Value = 4999.6665
$ws.Range("N31").Value = -5589.6665
$ws.Range("H34").Value = 2188.2273
$ws.Range("J34").Value = 4999.6665
$ws.Range("L34").Value = 4999.6665
$ws.Range("N34").Value = -5403.6665
$ws.Range("H58").Value = 2461.5386
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -1797
$ws.Range("N58").Value = -3906
$ws.Range("H122").Value = 6380.5386
$ws.Range("I122").Value = 4491.1665
$ws.Range("K122").Value = 13473.4995
$ws.Range("M122").Value = -11023.4995
$ws.Range("H136").Value = 2461.5386
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 169507.14
$ws.Range("I80").Value = 286152.5
$ws.Range("J80").Value = 13980
$ws.Range("K80").Value = 286152.5
$ws.Range("L80").Value = 13980
$ws.Range("M80").Value = -285154.5
$ws.Range("N80").Value = -15976
$ws.Range("H83").Value = 169507.14
$ws.Range("I83").Value = 286152.5
$ws.Range("J83").Value = 13980
$ws.Range("K83").Value = 1430762.5
$ws.Range("L83").Value = 69900
$ws.Range("M83").Value = -1425770.5
$ws.Range("N83").Value = -79884
$ws.Range("H102").Value = 3569
$ws.Range("I102").Value = 2799.8948
$ws.Range("K102").Value = 2799.8948
$ws.Range("M102").Value = -1177.8948
$ws.Range("H122").Value = 4256.8
$ws.Range("I122").Value = 2928
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 8784
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -6334
$ws.Range("N122").Value = -23650
$ws.Range("H126").Value = 3836.3333
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 4254.5
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 12763.5
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -17703.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2511.3333
$ws.Range("I40").Value = 2612.6
$ws.Range("J40").Value = 2005
$ws.Range("K40").Value = 2612.6
$ws.Range("L40").Value = 2005
$ws.Range("M40").Value = -2476.6
$ws.Range("N40").Value = -2277
$ws.Range("H95").Value = 58497.5
$ws.Range("J95").Value = 58497.5
$ws.Range("L95").Value = 58497.5
$ws.Range("N95").Value = -63989.5
$ws.Range("H98").Value = 200355
$ws.Range("J98").Value = 200355
$ws.Range("L98").Value = 200355
$ws.Range("N98").Value = -206345
$ws.Range("H100").Value = 54016.305
$ws.Range("I100").Value = 96864.664
$ws.Range("J100").Value = 7272.636
$ws.Range("K100").Value = 96864.664
$ws.Range("L100").Value = 7272.636
$ws.Range("M100").Value = -96323.664
$ws.Range("N100").Value = -8354.636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 34462.223
$ws.Range("I44").Value = 29499.5
$ws.Range("K44").Value = 29499.5
$ws.Range("M44").Value = -28945.5
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H101").Value = 19602
$ws.Range("J101").Value = 19602
$ws.Range("L101").Value = 19602
$ws.Range("N101").Value = -26092
$ws.Range("H103").Value = 53749.75
$ws.Range("J103").Value = 53749.75
$ws.Range("L103").Value = 53749.75
$ws.Range("N103").Value = -56093.75
$ws.Range("H122").Value = 1973.875
$ws.Range("I122").Value = 1798.6
$ws.Range("K122").Value = 5395.799999999999
$ws.Range("M122").Value = -2945.799999999999
